$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Fixes")
[void]$ws.Select()

# Fill in "Who" and "Which files are touched" for row 2 (the login/blank-password bug).
$ws.Range("C2").Value = "Ana Gorohovschi"
$ws.Range("D2").Value = "MarchMadnessGUI.java`nSerializeTest.java"
$ws.Range("D2").WrapText = $true

# Box every cell in the used range with a thin border (all borders).
$ws.Range("A1:D11").Borders.LineStyle = 1
$ws.Range("A1:D11").Borders.Weight = 2

# Leave the selection where the user ended up.
[void]$ws.Range("C11").Select()

Write-Host "done"
